$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a truly empty-but-text (shared-string) cell value,
# matching cells that hold an empty string rather than being blank.
function Set-EmptyTextCell($range) {
    $range.Formula = "'"
    $range.ClearFormats()
}

# Row 1: new header column I ("Other found locations")
$ws.Range("I1").Value = 'Other found locations'

# Row 2
$ws.Range("E2").Value = '[Tran Duc Anh%Ly%NULL%3,   Didier%Zanini%NULL%1,   Vincent%Laforge%NULL%1,   Sylvie%Arlotto%NULL%1,   Stephanie%Gentile%NULL%1,   Helene%Mendizabal%NULL%1,   Michael%Finaud%NULL%1,   David%Morel%NULL%1,   Olivier%Quenette%NULL%1,   Priscilla%Malfuson-Clot-Faybesse%NULL%1,   Alain%Midejean%NULL%1,   Phuc%Le-Dinh%NULL%1,   Gérard%Daher%NULL%1,   Berengere%Labarriere%NULL%1,   Anne-Marie%Morel-Roux%NULL%1,   Alain%Coquet%NULL%1,   Patrick%Augier%NULL%1,   Philippe%Parola%NULL%1,   Eric%Chabriere%NULL%1,   Didier%Raoult%NULL%1,   Philippe%Gautret%NULL%1]'
$ws.Range("I2").Value = '_PMC_elsevier'

# Row 3
$ws.Range("E3").Value = '[Matthieu%Mahévas%NULL%0,   Viet-Thi%Tran%NULL%2,   Viet-Thi%Tran%NULL%0,   Mathilde%Roumier%NULL%1,   Amélie%Chabrol%NULL%1,   Romain%Paule%NULL%1,   Constance%Guillaud%NULL%1,   Elena%Fois%NULL%1,   Raphael%Lepeule%NULL%1,   Tali-Anne%Szwebel%NULL%1,   François-Xavier%Lescure%NULL%1,   Frédéric%Schlemmer%NULL%1,   Marie%Matignon%NULL%1,   Mehdi%Khellaf%NULL%1,   Etienne%Crickx%NULL%1,   Benjamin%Terrier%NULL%1,   Caroline%Morbieu%NULL%1,   Paul%Legendre%NULL%1,   Julien%Dang%NULL%1,   Yoland%Schoindre%NULL%1,   Jean-Michel%Pawlotsky%NULL%1,   Marc%Michel%NULL%1,   Elodie%Perrodeau%NULL%1,   Nicolas%Carlier%NULL%1,   Nicolas%Roche%NULL%1,   Victoire%de Lastours%NULL%1,   Clément%Ourghanlian%NULL%1,   Solen%Kerneis%NULL%1,   Philippe%Ménager%NULL%1,   Luc%Mouthon%NULL%1,   Etienne%Audureau%NULL%1,   Philippe%Ravaud%NULL%1,   Bertrand%Godeau%NULL%1,   Sébastien%Gallien%NULL%1,   Nathalie%Costedoat-Chalumeau%NULL%1]'
$ws.Range("I3").Value = '_PMC'

# Row 4
$ws.Range("E4").Value = '[Ahmed%Hasan%coreGivesNoEmail%1,  Amir%Lerman%coreGivesNoEmail%1,  Charanjit%Rihal%coreGivesNoEmail%1,  Charles%Cagin%coreGivesNoEmail%1,  Derek%So%coreGivesNoEmail%1,  Erin%Iturriaga%coreGivesNoEmail%1,  Ivan%Chavez%coreGivesNoEmail%1,  J.%Dawn Abbott%coreGivesNoEmail%1,  Jang-Ho%Bae%coreGivesNoEmail%1,  Jean-Francois%Tanguay%coreGivesNoEmail%1,  Kent%Bailey%coreGivesNoEmail%1,  Liewei%Wang%coreGivesNoEmail%1,  Linnea%Baudhuin%coreGivesNoEmail%1,  Malcolm%Bell%coreGivesNoEmail%1,  Mandeep%Sidhu%coreGivesNoEmail%1,  Michael%E. Farkouh%coreGivesNoEmail%1,  Myung%Ho Jeong%coreGivesNoEmail%1,  Nancy%Geller%coreGivesNoEmail%1,  Naveen%L. Pereira%coreGivesNoEmail%1,  Paul%Gordon%coreGivesNoEmail%1,  Richard%Weinshilboum%coreGivesNoEmail%1,  Robert%Welsh%coreGivesNoEmail%1,  Ryan%Lennon%coreGivesNoEmail%1,  Shaun%G. Goodman%coreGivesNoEmail%1,  Verghese%Mathew%coreGivesNoEmail%1,  Yi-Ping%Fu%coreGivesNoEmail%1,  Yves%Rosenberg%coreGivesNoEmail%1]'
$ws.Range("F4").Value = 'not found'
$ws.Range("G4").Value = 'N/A'
Set-EmptyTextCell $ws.Range("I4")

# Row 5
$ws.Range("E5").Value = '[ Li%Yu%null%2,  Li%Yu%null%0]'
Set-EmptyTextCell $ws.Range("I5")

# Row 6
Set-EmptyTextCell $ws.Range("I6")

# Row 7
$ws.Range("E7").Value = '[Joseph%Magagnoli%NULL%1,   Siddharth%Narendran%NULL%1,   Felipe%Pereira%NULL%1,   Tammy H.%Cummings%NULL%1,   James W.%Hardin%NULL%1,   S. Scott%Sutton%NULL%1,   Jayakrishna%Ambati%NULL%1]'
$ws.Range("I7").Value = '_PMC_elsevier'

# Row 8
$ws.Range("E8").Value = '[ S.%Singh%null%1,    A.% Khan%null%1,    M.% Chowdhry%null%1,    A. % Chatterjee%null%1]'
$ws.Range("F8").Value = 'not found'
$ws.Range("G8").Value = 'N/A'
$ws.Range("I8").Value = '_MedBiorxiv'

